$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row above the footer row (117), pushing the footer row down to 118
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new day's data
$ws.Range("A117").Value = 43972
$ws.Range("B117").Value = 119
$ws.Range("C117").Value = 38465
$ws.Range("D117").Value = 36
$ws.Range("E117").Value = 7764

# Update the active selection to follow the footer row to its new position
[void]$ws.Range("B118").Select()

# Update the print area defined name to cover the new last row
$ws.PageSetup.PrintArea = "`$A`$1:`$E`$118"
